$p = $ppt.ActivePresentation

# Remove slides 3 through 10, keeping only slide 1 and slide 2.
for ($i = $p.Slides.Count; $i -ge 3; $i--) {
    $p.Slides.Item($i).Delete()
}

$s = $p.Slides.Item(2)

# Remove the picture and the title placeholder from slide 2.
$s.Shapes.Item("Picture 1").Delete()
$s.Shapes.Item("Title 1").Delete()
# Deleting a required layout placeholder leaves an auto-regenerated empty
# "ghost" placeholder behind (mirrors real PowerPoint behavior); delete it
# again so the title is actually gone.
$s.Shapes.Item("Title 3").Delete()

# Update the remaining content placeholder: single paragraph with the
# Quarto include directive, no bullet, zero indent/margin.
$shp = $s.Shapes.Item("Content Placeholder 2")
$shp.TextFrame.TextRange.Text = "{{ < include lci_nutrition_guts.qmd > }}"

$tr = $shp.TextFrame.TextRange
$tr.ParagraphFormat.Bullet.Type = 0

$lvl = $shp.TextFrame.Ruler.Levels.Item(1)
$lvl.LeftMargin = 0
$lvl.FirstMargin = 0
